$d = $word.ActiveDocument

# --- Change 3 (done first): remove the stray <w:bookmarkStart/><w:bookmarkEnd/>
# ("_GoBack") that sits alone in an empty paragraph further down the document,
# so the name is free to be re-created at its new location in change 2 below.
$d.Bookmarks("_GoBack").Delete()

# --- Change 1: remove <w:lang w:val="en-US"/> from the title run's rPr ---
$p1 = $d.Paragraphs.Item(1).Range
$xml1 = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidR="005E22B1" w:rsidRPr="00995179" w:rsidRDefault="00995179" w:rsidP="00995179">
<w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:u w:val="single"/></w:rPr></w:pPr>
<w:r w:rsidRPr="00995179"><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t>User Manual Kappa Kathy</w:t></w:r>
<w:r w:rsidRPr="00995179"><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> Universidad</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$p1.InsertXML($xml1) | Out-Null

# --- Change 2: split "preste" into "pres" + "en" + bookmark(_GoBack) + "te" ---
$p2 = $d.Paragraphs.Item(2).Range
$xml2 = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidR="00995179" w:rsidRDefault="008F149C" w:rsidP="00995179">
<w:r><w:t>El</w:t></w:r>
<w:r w:rsidR="00995179" w:rsidRPr="00995179"><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:t>pres</w:t></w:r>
<w:r><w:t>en</w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
<w:r><w:t>te</w:t></w:r>
<w:r w:rsidR="00995179" w:rsidRPr="00995179"><w:t xml:space="preserve"> modelo de datos fue </w:t></w:r>
<w:r w:rsidR="00995179"><w:t>creado</w:t></w:r>
<w:r w:rsidR="00995179" w:rsidRPr="00995179"><w:t xml:space="preserve"> con el propósito de almacenar información básica, y</w:t></w:r>
<w:r w:rsidR="00995179"><w:t xml:space="preserve"> necesaria para los estudiantes de la universidad. Dicha información será tomada desde el momento de la inscripción hasta la toma de las materias respectivas.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$p2.InsertXML($xml2) | Out-Null

# --- Change 4: merge the "dd" / proofErr-wrapped run back into a single run ---
$p26 = $d.Paragraphs.Item(26).Range
$xml4 = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidR="0072579A" w:rsidRDefault="0072579A" w:rsidP="0072579A">
<w:pPr>
<w:pStyle w:val="Prrafodelista"/>
<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
<w:tabs><w:tab w:val="left" w:pos="8436"/></w:tabs>
</w:pPr>
<w:r><w:t>Si desea ingresar una fecha se le aconseja respetar el formato definido para tener un funcionamiento más efectivo (dd/mm/yyyy)</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$p26.InsertXML($xml4) | Out-Null

Write-Host "edits applied"
